$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178-186 down to 179-187
$ws.Rows.Item(178).Insert()

# Copy the fixed/common column values from the row above (now row 179, previously row 178)
# into the new row 178
$ws.Range("A178").Value = 11
$ws.Range("B178").Value = "Vega Monumental Concepción"
$ws.Range("C178").Value = "Bíobío"
$ws.Range("D178").Value = 45008
$ws.Range("D178").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100102
$ws.Range("H178").Value = "Cítricos"
$ws.Range("I178").Value = 100102004
$ws.Range("J178").Value = "Mandarina"
$ws.Range("K178").Value = "Murcott"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 370
$ws.Range("N178").Value = 9000
$ws.Range("O178").Value = 10000
$ws.Range("P178").Value = 9405
$ws.Range("Q178").Value = "`$/caja 15 kilos granel"
$ws.Range("R178").Value = "Región de O'Higgins"
$ws.Range("S178").Value = 627
$ws.Range("T178").Value = 15
